$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6316415812111416, 0.5980219932806392, 0.5775950176231675, 0.5693256787541543, 0.5679558902484416, 0.5774832716569165, 0.620005155105531, 0.7050789492779757, 0.7685892072335037, 0.7976958320494987, 0.8087482222837252, 0.8063665530126798, 0.798604514407657, 0.7938539741218165, 0.7666913019581614, 0.7500826156762344, 0.7405500599282107, 0.737326010204157, 0.7518485401755868, 0.800883594921487, 0.8331074030706986, 0.8158930068223071, 0.7510501158257057, 0.6818858339525207)
    3 = @(0.2499689866440526, 0.2507483669695638, 0.251256948581867, 0.251471773769115, 0.2515079033874663, 0.2512598150944072, 0.2502314947151447, 0.2484523869822048, 0.247288752190741, 0.2467902728652476, 0.2466059296246073, 0.2466454349421063, 0.2467750183535387, 0.246854967023971, 0.2473219484925018, 0.2476163185242442, 0.2477885385148504, 0.2478473489298167, 0.2475846816890908, 0.2467368367491787, 0.2462084763161414, 0.2464881215504775, 0.2475989754119503, 0.2489083964042145)
    5 = @(0.2247082988006497, 0.2246401704352969, 0.2246703141947606, 0.224700748036053, 0.2247068992788002, 0.2246706510789345, 0.2246699001405119, 0.2252374330408209, 0.2259987364270692, 0.2264193207178842, 0.2265892160291685, 0.2265521540918876, 0.2264330853682743, 0.2263615350905148, 0.2259727403650551, 0.2257532125868416, 0.225633939028846, 0.2255947576279169, 0.2257758583044556, 0.2264677706629996, 0.2269819072371497, 0.226701852551674, 0.2257655985675804, 0.2250232477095722)
    6 = @(2.482278692466352, 2.486648447878025, 2.490421521384903, 2.49223350082989, 2.492550962583962, 2.49044484682107, 2.48355934016223, 2.478693475755989, 2.480370120293784, 2.482270545451826, 2.483153514304519, 2.4829560920072, 2.482339916292887, 2.481983751450272, 2.4802687925729, 2.479507915253023, 2.479177382419678, 2.479083873438015, 2.479577829156341, 2.482516471857153, 2.485388861871598, 2.483768815406364, 2.479545888055114, 2.47908660878862)
    7 = @(0.002511640963699269, 0.00251412383306227, 0.002515729984375064, 0.002516405098437591, 0.002516518446219159, 0.002515739005747919, 0.002512480145121716, 0.002506734700793566, 0.002502902965665231, 0.002501243550332621, 0.002500627141204735, 0.002500759364160166, 0.002501192598270264, 0.002501459524496706, 0.002503013091212808, 0.002503987541912695, 0.002504555898100197, 0.002504749688405014, 0.002503882994977008, 0.00250106502226294, 0.002499293090582045, 0.002500232437423657, 0.002503930235149757, 0.002508220327561111)
    10 = @(0.1016191822282124, 0.1004741009542549, 0.09976624546303015, 0.09947660648433754, 0.09942844127355954, 0.09976234405344186, 0.1012253586216012, 0.1040558267121305, 0.1061112259522794, 0.107040905753415, 0.1073921692607911, 0.1073165536476424, 0.1070698202784541, 0.106918586023113, 0.1060503602478704, 0.1055163524983094, 0.10520870468185, 0.1051044550159759, 0.1055732504939684, 0.1071423133370786, 0.1081631959107128, 0.1076187581077832, 0.1055475289007504, 0.1032943017301733)
    11 = @(0.2619705179571525, 0.2337415419360411, 0.2164239432800343, 0.2093710843505363, 0.2082002283000719, 0.2163288084019541, 0.2522342649839118, 0.3227503885056251, 0.3746090607895383, 0.3982092593583673, 0.4071470610232097, 0.4052221116022281, 0.3989445617675074, 0.3950994879228631, 0.3730668846302194, 0.3595526994470788, 0.3517806110440631, 0.3491492856905722, 0.3609912161607554, 0.4007884088674984, 0.42680332781066, 0.4129183484985219, 0.360340870832033, 0.3036640457974329)
    12 = @(0.3666456695819562, 0.3619205781731125, 0.3591531091862095, 0.3580591296170326, 0.3578795202534337, 0.3591382183920686, 0.3649887861668901, 0.3775172160872131, 0.3873587979172584, 0.3919729983736602, 0.3937398716669662, 0.3933584757314463, 0.3921179686675345, 0.3913606669909626, 0.3870599968786337, 0.3844567083979769, 0.3829722903898158, 0.3824719165298944, 0.384732496375193, 0.3924818056227082, 0.3976604573597768, 0.3948861310900327, 0.384607774471263, 0.374015620599252)
    13 = @(0.2302845244533209, 0.2228534681675818, 0.2183778133635492, 0.2165759947571146, 0.2162781404313279, 0.2183534239283098, 0.2277043229524054, 0.2467261818779534, 0.2611127396596373, 0.2677456231846378, 0.2702698926457856, 0.2697256909617209, 0.2679530459763058, 0.2668688779353658, 0.2606810290328028, 0.256907498262585, 0.2547453905509798, 0.2540147733457019, 0.2573083364100626, 0.2684733757601236, 0.2758433925478698, 0.2719032496531426, 0.2571270945913682, 0.2415076037814003)
    15 = @(5.049977378944504, 5.082467025354518, 5.104616984303476, 5.114197030583753, 5.115821243049751, 5.104743941683353, 5.060723102516647, 4.991854050425928, 4.951886344891335, 4.936010156161331, 4.930329576054334, 4.93153825386139, 4.935536170561164, 4.938028163093804, 4.952970270325068, 4.962727135592786, 4.968556020928048, 4.970566853867439, 4.961666043336066, 4.934352892648064, 4.918433700062366, 4.92675337521564, 4.962145079294885, 5.008617351203725)
}

foreach ($colIdx in $data.Keys) {
    $vals = $data[$colIdx]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($i + 2, [int]$colIdx).Value = $vals[$i]
    }
}
